$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"12.44313933333333"
$ws.Range("H2").Value = [double]"37.329418"
$ws.Range("I2").Value = [double]"0.9304541596872169"
$ws.Range("J2").Value = [double]"0.930454159687217"
$ws.Range("M2").Value = [double]"0.7861523333333333"
$ws.Range("N2").Value = [double]"2.358457"
$ws.Range("O2").Value = [double]"0.01668797875153133"
$ws.Range("P2").Value = [double]"0.01668797875153133"
$ws.Range("Q2").Value = [double]"9.782203020891778"
$ws.Range("R2").Value = [double]"88.039827188026"
$ws.Range("S2").Value = [double]"0.01552739924613421"
$ws.Range("T2").Value = [double]"0.01552739924613421"
$ws.Range("G3").Value = [double]"12.44313933333333"
$ws.Range("H3").Value = [double]"37.329418"
$ws.Range("I3").Value = [double]"0.9304541596872169"
$ws.Range("J3").Value = [double]"0.930454159687217"
$ws.Range("O3").Value = [double]"0.5736784050900728"
$ws.Range("P3").Value = [double]"0.5736784050900727"
$ws.Range("Q3").Value = [double]"336.280307570354"
$ws.Range("R3").Value = [double]"3026.522768133186"
$ws.Range("S3").Value = [double]"0.5337814583387865"
$ws.Range("T3").Value = [double]"0.5337814583387865"
$ws.Range("G4").Value = [double]"12.44313933333333"
$ws.Range("H4").Value = [double]"37.329418"
$ws.Range("I4").Value = [double]"0.9304541596872169"
$ws.Range("J4").Value = [double]"0.930454159687217"
$ws.Range("M4").Value = [double]"19.22475933333333"
$ws.Range("N4").Value = [double]"57.674278"
$ws.Range("O4").Value = [double]"0.4080918692916219"
$ws.Range("P4").Value = [double]"0.4080918692916219"
$ws.Range("Q4").Value = [double]"239.2163590344671"
$ws.Range("R4").Value = [double]"2152.947231310204"
$ws.Range("S4").Value = [double]"0.3797107773169217"
$ws.Range("T4").Value = [double]"0.3797107773169217"
$ws.Range("G5").Value = [double]"12.44313933333333"
$ws.Range("H5").Value = [double]"37.329418"
$ws.Range("I5").Value = [double]"0.9304541596872169"
$ws.Range("J5").Value = [double]"0.930454159687217"
$ws.Range("M5").Value = [double]"0.07263"
$ws.Range("N5").Value = [double]"0.21789"
$ws.Range("O5").Value = [double]"0.00154174686677398"
$ws.Range("P5").Value = [double]"0.00154174686677398"
$ws.Range("Q5").Value = [double]"0.9037452097799999"
$ws.Range("R5").Value = [double]"8.133706888019999"
$ws.Range("S5").Value = [double]"0.001434524785374583"
$ws.Range("T5").Value = [double]"0.001434524785374583"
$ws.Range("G6").Value = [double]"0.8272856666666667"
$ws.Range("I6").Value = [double]"0.0618615101204856"
$ws.Range("J6").Value = [double]"0.06186151012048561"
$ws.Range("M6").Value = [double]"0.7861523333333333"
$ws.Range("N6").Value = [double]"2.358457"
$ws.Range("O6").Value = [double]"0.01668797875153133"
$ws.Range("P6").Value = [double]"0.01668797875153133"
$ws.Range("Q6").Value = [double]"0.6503725571832223"
$ws.Range("R6").Value = [double]"5.853353014649"
$ws.Range("S6").Value = [double]"0.001032343566428304"
$ws.Range("T6").Value = [double]"0.001032343566428304"
$ws.Range("G7").Value = [double]"0.8272856666666667"
$ws.Range("I7").Value = [double]"0.0618615101204856"
$ws.Range("J7").Value = [double]"0.06186151012048561"
$ws.Range("O7").Value = [double]"0.5736784050900728"
$ws.Range("P7").Value = [double]"0.5736784050900727"
$ws.Range("S7").Value = [double]"0.03548861246238358"
$ws.Range("T7").Value = [double]"0.03548861246238358"
$ws.Range("G8").Value = [double]"0.8272856666666667"
$ws.Range("I8").Value = [double]"0.0618615101204856"
$ws.Range("J8").Value = [double]"0.06186151012048561"
$ws.Range("M8").Value = [double]"19.22475933333333"
$ws.Range("N8").Value = [double]"57.674278"
$ws.Range("O8").Value = [double]"0.4080918692916219"
$ws.Range("P8").Value = [double]"0.4080918692916219"
$ws.Range("Q8").Value = [double]"15.90436784158289"
$ws.Range("R8").Value = [double]"143.139310574246"
$ws.Range("S8").Value = [double]"0.02524517930227156"
$ws.Range("T8").Value = [double]"0.02524517930227156"
$ws.Range("G9").Value = [double]"0.8272856666666667"
$ws.Range("I9").Value = [double]"0.0618615101204856"
$ws.Range("J9").Value = [double]"0.06186151012048561"
$ws.Range("M9").Value = [double]"0.07263"
$ws.Range("N9").Value = [double]"0.21789"
$ws.Range("O9").Value = [double]"0.00154174686677398"
$ws.Range("P9").Value = [double]"0.00154174686677398"
$ws.Range("Q9").Value = [double]"0.06008575797"
$ws.Range("R9").Value = [double]"0.54077182173"
$ws.Range("S9").Value = [double]"9.537478940216555E-05"
$ws.Range("T9").Value = [double]"9.537478940216555E-05"
$ws.Range("G10").Value = [double]"0.073169"
$ws.Range("H10").Value = [double]"0.219507"
$ws.Range("I10").Value = [double]"0.005471320266243153"
$ws.Range("J10").Value = [double]"0.005471320266243153"
$ws.Range("M10").Value = [double]"0.7861523333333333"
$ws.Range("N10").Value = [double]"2.358457"
$ws.Range("O10").Value = [double]"0.01668797875153133"
$ws.Range("P10").Value = [double]"0.01668797875153133"
$ws.Range("Q10").Value = [double]"0.05752198007766667"
$ws.Range("R10").Value = [double]"0.5176978206990001"
$ws.Range("S10").Value = [double]"9.130527634588847E-05"
$ws.Range("T10").Value = [double]"9.130527634588845E-05"
$ws.Range("G11").Value = [double]"0.073169"
$ws.Range("H11").Value = [double]"0.219507"
$ws.Range("I11").Value = [double]"0.005471320266243153"
$ws.Range("J11").Value = [double]"0.005471320266243153"
$ws.Range("O11").Value = [double]"0.5736784050900728"
$ws.Range("P11").Value = [double]"0.5736784050900727"
$ws.Range("Q11").Value = [double]"1.977418492671"
$ws.Range("R11").Value = [double]"17.796766434039"
$ws.Range("S11").Value = [double]"0.003138778284075364"
$ws.Range("T11").Value = [double]"0.003138778284075364"
$ws.Range("G12").Value = [double]"0.073169"
$ws.Range("H12").Value = [double]"0.219507"
$ws.Range("I12").Value = [double]"0.005471320266243153"
$ws.Range("J12").Value = [double]"0.005471320266243153"
$ws.Range("M12").Value = [double]"19.22475933333333"
$ws.Range("N12").Value = [double]"57.674278"
$ws.Range("O12").Value = [double]"0.4080918692916219"
$ws.Range("P12").Value = [double]"0.4080918692916219"
$ws.Range("Q12").Value = [double]"1.406656415660667"
$ws.Range("R12").Value = [double]"12.659907740946"
$ws.Range("S12").Value = [double]"0.002232801314944303"
$ws.Range("T12").Value = [double]"0.002232801314944303"
$ws.Range("G13").Value = [double]"0.073169"
$ws.Range("H13").Value = [double]"0.219507"
$ws.Range("I13").Value = [double]"0.005471320266243153"
$ws.Range("J13").Value = [double]"0.005471320266243153"
$ws.Range("M13").Value = [double]"0.07263"
$ws.Range("N13").Value = [double]"0.21789"
$ws.Range("O13").Value = [double]"0.00154174686677398"
$ws.Range("P13").Value = [double]"0.00154174686677398"
$ws.Range("Q13").Value = [double]"0.00531426447"
$ws.Range("R13").Value = [double]"0.04782838023"
$ws.Range("S13").Value = [double]"8.43539087759736E-06"
$ws.Range("T13").Value = [double]"8.43539087759736E-06"
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.029595"
$ws.Range("H14").Value = [double]"0.088785"
$ws.Range("I14").Value = [double]"0.002213009926054287"
$ws.Range("J14").Value = [double]"0.002213009926054287"
$ws.Range("M14").Value = [double]"0.7861523333333333"
$ws.Range("N14").Value = [double]"2.358457"
$ws.Range("O14").Value = [double]"0.01668797875153133"
$ws.Range("P14").Value = [double]"0.01668797875153133"
$ws.Range("Q14").Value = [double]"0.023266178305"
$ws.Range("R14").Value = [double]"0.209395604745"
$ws.Range("S14").Value = [double]"3.693066262292186E-05"
$ws.Range("T14").Value = [double]"3.693066262292185E-05"
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.029595"
$ws.Range("H15").Value = [double]"0.088785"
$ws.Range("I15").Value = [double]"0.002213009926054287"
$ws.Range("J15").Value = [double]"0.002213009926054287"
$ws.Range("O15").Value = [double]"0.5736784050900728"
$ws.Range("P15").Value = [double]"0.5736784050900727"
$ws.Range("Q15").Value = [double]"0.799815499605"
$ws.Range("R15").Value = [double]"7.198339496445"
$ws.Range("S15").Value = [double]"0.001269556004827323"
$ws.Range("T15").Value = [double]"0.001269556004827323"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.029595"
$ws.Range("H16").Value = [double]"0.088785"
$ws.Range("I16").Value = [double]"0.002213009926054287"
$ws.Range("J16").Value = [double]"0.002213009926054287"
$ws.Range("M16").Value = [double]"19.22475933333333"
$ws.Range("N16").Value = [double]"57.674278"
$ws.Range("O16").Value = [double]"0.4080918692916219"
$ws.Range("P16").Value = [double]"0.4080918692916219"
$ws.Range("Q16").Value = [double]"0.56895675247"
$ws.Range("R16").Value = [double]"5.12061077223"
$ws.Range("S16").Value = [double]"0.0009031113574844079"
$ws.Range("T16").Value = [double]"0.0009031113574844078"
$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.3333333333333333"
$ws.Range("G17").Value = [double]"0.029595"
$ws.Range("H17").Value = [double]"0.088785"
$ws.Range("I17").Value = [double]"0.002213009926054287"
$ws.Range("J17").Value = [double]"0.002213009926054287"
$ws.Range("M17").Value = [double]"0.07263"
$ws.Range("N17").Value = [double]"0.21789"
$ws.Range("O17").Value = [double]"0.00154174686677398"
$ws.Range("P17").Value = [double]"0.00154174686677398"
$ws.Range("Q17").Value = [double]"0.00214948485"
$ws.Range("R17").Value = [double]"0.01934536365"
$ws.Range("S17").Value = [double]"3.411901119633915E-06"
$ws.Range("T17").Value = [double]"3.411901119633914E-06"
